$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended after the last existing row (09/15/2025 -> 09/16/2025).
# Force column A to text first so the MM/DD/YYYY-looking date string is kept
# as a literal string (matching the existing rows) instead of being
# auto-converted into a date serial number; then drop the formatting again
# so the new row's cells carry no explicit style, same as rows 2-14.
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "09/16/2025"
$ws.Range("A15").ClearFormats()

$ws.Range("B15").Value = 0.1244621484369293
$ws.Range("C15").Value = 0.8755378515630707
